$d = $word.ActiveDocument

# The resume's header line ("Dheeraj Chand") is immediately followed by the
# PROFESSIONAL SUMMARY heading. The short-resume contact line was dropped
# during generation, so add it back as its own centered paragraph right
# after the name, mirroring the long-resume layout.
#
# Using Find/Replace with a literal "^p" (paragraph mark) in the replacement
# text inserts a brand-new paragraph without carrying over the name run's
# bold/28pt character formatting - the new run ends up with the document's
# default (unformatted) run properties, and the new paragraph inherits the
# centered alignment already set on that paragraph.
$find = $d.Content.Find
[void]$find.Execute(
    "Dheeraj Chand",
    $false,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "Dheeraj Chand^p202.550.7110 | dheeraj.chand@gmail.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/ | Austin, TX",
    2
)
